$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Capture the text currently sitting in D3/D4 before we touch the merges,
# since those are the values that must move into C3/C4.
$appNameValue = $ws.Range("D3").Text
$snapshotDateValue = $ws.Range("D4").Text

# Un-merge the label/value header blocks in row 3 and row 4 (F1:G1 stays merged).
$ws.Range("B3:C3").UnMerge()
$ws.Range("B4:C4").UnMerge()

# Slide the RepGen placeholder text one column to the left (D -> C) and
# clear out the now-empty D cells while keeping their formatting.
$ws.Range("C3").Value = $appNameValue
$ws.Range("D3").ClearContents()

$ws.Range("C4").Value = $snapshotDateValue
$ws.Range("D4").ClearContents()

# Move the active selection to B7.
$ws.Range("B7").Select()
